# Apply the "scannertimes" edits described in the commit diff.
#
# Changes:
#  - A1 (roomwidth):  12 -> 14
#  - A2 (roomlength):  16 -> 54
#  - Column B/C widths widened (bestFit removed, explicit customWidth set)
#  - D7 gets a new label "bytes/day at 4cm" (adds a 4th shared string)
#  - Selected cell moves from C10 to C9
#  - All dependent formulas (B1, B2, B4, C4, C5, C7) recalculate automatically

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the two input values that drive the formulas ---
$ws.Range("A1").Value = 14
$ws.Range("A2").Value = 54

# --- Label the new "bytes/day at 4cm" total in D7 ---
$ws.Range("D7").Value = "bytes/day at 4cm"

# --- Resize columns B and C (removes bestFit, sets explicit customWidth) ---
# The target stored OOXML widths are 15.5703125 and 17.140625 characters.
# Excel's COM ColumnWidth setter quantizes to whole pixels (1/6-character
# steps for the default Calibri 11 / MDW=7 grid), so we pass the
# ColumnWidth value whose pixel-quantized result lands as close as
# possible to those targets (15.5 and 17.1666... respectively).
$ws.Range("B1").ColumnWidth = 14.666666666666666
$ws.Range("C1").ColumnWidth = 16.333333333333332

# --- Move the active selection from C10 to C9 ---
$ws.Range("C9").Select() | Out-Null
